# Weekly update: insert 3 new rows of "Pimiento" price data (week of 2021-09-22,
# serial 44461) above the existing data block, pushing rows 466:487 down to 469:490.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 466:468 - everything currently at 466 and below
# (through 487) shifts down to 469:490.
$ws.Range("A466:A468").EntireRow.Insert()

# Populate the three newly inserted rows with the new week's data.
$newRows = @(
    @{ Row = 466; H = "Cuatro cascos verde"; I = "Primera"; J = 800;  K = 30000; L = 31000; M = 30500; P = 1694 },
    @{ Row = 467; H = "Cuatro cascos verde"; I = "Segunda"; J = 600;  K = 28000; L = 29000; M = 28500; P = 1583 },
    @{ Row = 468; H = "Cuatro cascos verde"; I = "Tercera"; J = 500;  K = 25000; L = 26000; M = 25500; P = 1417 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44461
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = 100112002
    $ws.Cells.Item($row, 7).Value = "Pimiento"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "$/caja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
